# "updated main GSC export data"
# Append two new daily rows (2025-12-02, 2025-12-03) to the bottom of the
# "Chart" sheet's date/URL-count table, following the same pattern as every
# preceding row (B = 0 Non-HTTPS URLs, C = 27 HTTPS URLs).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

$newDates = "2025-12-02", "2025-12-03"
$startRow = 58

for ($i = 0; $i -lt $newDates.Length; $i++) {
    $row = $startRow + $i

    # Pre-format column A as Text so the date-like string ("yyyy-MM-dd")
    # is stored as a literal string (matching every other row in the
    # column) instead of being auto-parsed into a date serial number.
    $ws.Range("A$row").NumberFormat = "@"
    $ws.Range("A$row").Value = $newDates[$i]
    # Drop back to the sheet's default (General) formatting now that the
    # text value is locked in, same as the untouched data rows above it.
    $ws.Range("A$row").ClearFormats()

    $ws.Range("B$row").Value = 0
    $ws.Range("C$row").Value = 27
}
